# Applies the referee-stats update described in the commit:
#   chore(runtime): publish files + archive (2025-12-04 19:17:09)
#
# Updates per-referee aggregate stats (games played, PIM totals, penalty
# breakdowns, etc.) on the "Главные" and "Линейные" sheets, and refreshes the
# as_of_utc timestamp (column AA) for every data row on both sheets.

$wb  = $excel.ActiveWorkbook
$ws2 = $wb.Worksheets.Item(2)   # Главные
$ws3 = $wb.Worksheets.Item(3)   # Линейные

function Set-Cell {
    param($ws, [string]$ref, $value)
    $ws.Range($ref).Value = $value
}

# ---- Главные (Worksheets.Item(2)): updated stat cells ----
Set-Cell $ws2 "C2" 33
Set-Cell $ws2 "D2" 767
Set-Cell $ws2 "E2" 339
Set-Cell $ws2 "F2" 428
Set-Cell $ws2 "G2" 23.24
Set-Cell $ws2 "H2" 10.27
Set-Cell $ws2 "I2" 12.97
Set-Cell $ws2 "J2" 132
Set-Cell $ws2 "K2" 159
Set-Cell $ws2 "L2" 9
Set-Cell $ws2 "M2" 6
Set-Cell $ws2 "O2" 2
Set-Cell $ws2 "P2" 1
Set-Cell $ws2 "V2" 14
Set-Cell $ws2 "W2" 18
Set-Cell $ws2 "C4" 23
Set-Cell $ws2 "D4" 354
Set-Cell $ws2 "E4" 151
Set-Cell $ws2 "F4" 203
Set-Cell $ws2 "G4" 15.39
Set-Cell $ws2 "H4" 6.57
Set-Cell $ws2 "I4" 8.83
Set-Cell $ws2 "J4" 73
Set-Cell $ws2 "K4" 89
Set-Cell $ws2 "C9" 32
Set-Cell $ws2 "D9" 494
Set-Cell $ws2 "E9" 259
Set-Cell $ws2 "F9" 235
Set-Cell $ws2 "G9" 15.44
Set-Cell $ws2 "H9" 8.09
Set-Cell $ws2 "I9" 7.34
Set-Cell $ws2 "J9" 127
Set-Cell $ws2 "K9" 115
Set-Cell $ws2 "C12" 21
Set-Cell $ws2 "D12" 349
Set-Cell $ws2 "E12" 151
Set-Cell $ws2 "F12" 198
Set-Cell $ws2 "G12" 16.62
Set-Cell $ws2 "H12" 7.19
Set-Cell $ws2 "I12" 9.43
Set-Cell $ws2 "J12" 63
Set-Cell $ws2 "K12" 74
Set-Cell $ws2 "C15" 21
Set-Cell $ws2 "D15" 376
Set-Cell $ws2 "E15" 179
Set-Cell $ws2 "F15" 197
Set-Cell $ws2 "H15" 8.52
Set-Cell $ws2 "I15" 9.380000000000001
Set-Cell $ws2 "J15" 67
Set-Cell $ws2 "K15" 86
Set-Cell $ws2 "C17" 20
Set-Cell $ws2 "D17" 302
Set-Cell $ws2 "E17" 113
Set-Cell $ws2 "F17" 189
Set-Cell $ws2 "G17" 15.1
Set-Cell $ws2 "H17" 5.65
Set-Cell $ws2 "I17" 9.449999999999999
Set-Cell $ws2 "J17" 54
Set-Cell $ws2 "K17" 77
Set-Cell $ws2 "C20" 31
Set-Cell $ws2 "D20" 523
Set-Cell $ws2 "E20" 230
Set-Cell $ws2 "F20" 293
Set-Cell $ws2 "G20" 16.87
Set-Cell $ws2 "H20" 7.42
Set-Cell $ws2 "I20" 9.449999999999999
Set-Cell $ws2 "J20" 105
Set-Cell $ws2 "K20" 109
Set-Cell $ws2 "C24" 31
Set-Cell $ws2 "D24" 573
Set-Cell $ws2 "E24" 274
Set-Cell $ws2 "F24" 299
Set-Cell $ws2 "G24" 18.48
Set-Cell $ws2 "H24" 8.84
Set-Cell $ws2 "I24" 9.65
Set-Cell $ws2 "J24" 112
Set-Cell $ws2 "K24" 122
Set-Cell $ws2 "L24" 6
Set-Cell $ws2 "M24" 5
Set-Cell $ws2 "O24" 1
Set-Cell $ws2 "P24" 1
Set-Cell $ws2 "V24" 20
Set-Cell $ws2 "W24" 14
Set-Cell $ws2 "C25" 31
Set-Cell $ws2 "D25" 479
Set-Cell $ws2 "E25" 232
Set-Cell $ws2 "F25" 247
Set-Cell $ws2 "G25" 15.45
Set-Cell $ws2 "H25" 7.48
Set-Cell $ws2 "I25" 7.97
Set-Cell $ws2 "J25" 111
Set-Cell $ws2 "K25" 116
Set-Cell $ws2 "C26" 12
Set-Cell $ws2 "D26" 256
Set-Cell $ws2 "E26" 114
Set-Cell $ws2 "F26" 142
Set-Cell $ws2 "G26" 21.33
Set-Cell $ws2 "H26" 9.5
Set-Cell $ws2 "I26" 11.83
Set-Cell $ws2 "J26" 52
Set-Cell $ws2 "K26" 51

# ---- Линейные (Worksheets.Item(3)): updated stat cells ----
Set-Cell $ws3 "C4" 15
Set-Cell $ws3 "D4" 224
Set-Cell $ws3 "E4" 102
Set-Cell $ws3 "F4" 122
Set-Cell $ws3 "G4" 14.93
Set-Cell $ws3 "H4" 6.8
Set-Cell $ws3 "I4" 8.130000000000001
Set-Cell $ws3 "J4" 51
Set-Cell $ws3 "K4" 51
Set-Cell $ws3 "C7" 18
Set-Cell $ws3 "D7" 349
Set-Cell $ws3 "E7" 144
Set-Cell $ws3 "F7" 205
Set-Cell $ws3 "G7" 19.39
Set-Cell $ws3 "H7" 8
Set-Cell $ws3 "I7" 11.39
Set-Cell $ws3 "J7" 52
Set-Cell $ws3 "K7" 65
Set-Cell $ws3 "L7" 4
Set-Cell $ws3 "M7" 3
Set-Cell $ws3 "O7" 2
Set-Cell $ws3 "P7" 1
Set-Cell $ws3 "V7" 12
Set-Cell $ws3 "W7" 10
Set-Cell $ws3 "C14" 29
Set-Cell $ws3 "D14" 454
Set-Cell $ws3 "E14" 231
Set-Cell $ws3 "F14" 223
Set-Cell $ws3 "G14" 15.66
Set-Cell $ws3 "H14" 7.97
Set-Cell $ws3 "I14" 7.69
Set-Cell $ws3 "J14" 113
Set-Cell $ws3 "K14" 104
Set-Cell $ws3 "C15" 26
Set-Cell $ws3 "D15" 505
Set-Cell $ws3 "E15" 259
Set-Cell $ws3 "F15" 246
Set-Cell $ws3 "G15" 19.42
Set-Cell $ws3 "H15" 9.960000000000001
Set-Cell $ws3 "I15" 9.460000000000001
Set-Cell $ws3 "J15" 107
Set-Cell $ws3 "K15" 103
Set-Cell $ws3 "C16" 30
Set-Cell $ws3 "D16" 484
Set-Cell $ws3 "E16" 224
Set-Cell $ws3 "F16" 260
Set-Cell $ws3 "G16" 16.13
Set-Cell $ws3 "H16" 7.47
Set-Cell $ws3 "I16" 8.67
Set-Cell $ws3 "J16" 102
Set-Cell $ws3 "K16" 110
Set-Cell $ws3 "C19" 29
Set-Cell $ws3 "D19" 465
Set-Cell $ws3 "E19" 224
Set-Cell $ws3 "F19" 241
Set-Cell $ws3 "G19" 16.03
Set-Cell $ws3 "H19" 7.72
Set-Cell $ws3 "I19" 8.31
Set-Cell $ws3 "J19" 107
Set-Cell $ws3 "K19" 108
Set-Cell $ws3 "C21" 32
Set-Cell $ws3 "D21" 607
Set-Cell $ws3 "E21" 252
Set-Cell $ws3 "F21" 355
Set-Cell $ws3 "G21" 18.97
Set-Cell $ws3 "H21" 7.88
Set-Cell $ws3 "I21" 11.09
Set-Cell $ws3 "J21" 116
Set-Cell $ws3 "K21" 145
Set-Cell $ws3 "C23" 17
Set-Cell $ws3 "D23" 262
Set-Cell $ws3 "E23" 128
Set-Cell $ws3 "F23" 134
Set-Cell $ws3 "G23" 15.41
Set-Cell $ws3 "H23" 7.53
Set-Cell $ws3 "I23" 7.88
Set-Cell $ws3 "J23" 59
Set-Cell $ws3 "K23" 62
Set-Cell $ws3 "C24" 32
Set-Cell $ws3 "D24" 550
Set-Cell $ws3 "E24" 219
Set-Cell $ws3 "F24" 331
Set-Cell $ws3 "G24" 17.19
Set-Cell $ws3 "H24" 6.84
Set-Cell $ws3 "I24" 10.34
Set-Cell $ws3 "J24" 97
Set-Cell $ws3 "K24" 128
Set-Cell $ws3 "C26" 29
Set-Cell $ws3 "D26" 572
Set-Cell $ws3 "E26" 248
Set-Cell $ws3 "F26" 324
Set-Cell $ws3 "G26" 19.72
Set-Cell $ws3 "H26" 8.550000000000001
Set-Cell $ws3 "I26" 11.17
Set-Cell $ws3 "J26" 99
Set-Cell $ws3 "K26" 102

# ---- refresh as_of_utc (column AA) for every data row (2-26) on both sheets ----
$newTimestamp = "2025-12-04 11:17:07"
for ($r = 2; $r -le 26; $r++) {
    Set-Cell $ws2 "AA$r" $newTimestamp
    Set-Cell $ws3 "AA$r" $newTimestamp
}

"Updated referee stats and timestamps on Главные/Линейные sheets."